$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Marking row (B11): correct marks value changed from 3 to 5
$ws.Range("B11").Value = 5

# Update Total row (B12): total marks changed from 51 to 85
$ws.Range("B12").Value = 85

# Update Total row (E12): corr/total marks text changed from "42/84" to "85/140"
$ws.Range("E12").Value = "85/140"
